$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 & 9: swap the oro_type B/C values (policy moments reordering)
$ws.Range("B8").Value = "Conservation"
$ws.Range("C8").Value = "Conservation"
$ws.Range("B9").Value = "Human assisted evolution"
$ws.Range("C9").Value = "Human-assisted evolution"

# E5:E11 - replace formulas with distinct literal colour codes for each ORO type
$ws.Range("E5").Value = "#026996"
$ws.Range("E6").Value = "#0688c2"
$ws.Range("E7").Value = "#9ed7f0"
$ws.Range("E8").Value = "#078257"
$ws.Range("E9").Value = "#43b08a"
$ws.Range("E10").Value = "#600787"
$ws.Range("E11").Value = "#ad5ad1"

# Update sheet view: scroll back to top and change selection to D10
$ws.Range("D10").Select()
